$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5899.7144
$ws.Range("J40").Value = 6939.6
$ws.Range("L40").Value = 6939.6
$ws.Range("N40").Value = -7289.6

$ws.Range("H107").Value = 783.8333
$ws.Range("I107").Value = 650.625
$ws.Range("K107").Value = 650.625
$ws.Range("M107").Value = 1269.375

$ws.Range("H112").Value = 2389.7273
$ws.Range("J112").Value = 2999.1428
$ws.Range("L112").Value = 8997.428400000001
$ws.Range("N112").Value = -11213.4284

$ws.Range("H125").Value = 35717236
$ws.Range("I125").Value = 41668444
$ws.Range("J125").Value = 10000
$ws.Range("K125").Value = 375015996
$ws.Range("L125").Value = 90000
$ws.Range("M125").Value = -375013536
$ws.Range("N125").Value = -94920

$ws.Range("H137").Value = 2716.5833
$ws.Range("I137").Value = 999.75
$ws.Range("K137").Value = 2999.25
$ws.Range("M137").Value = -449.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 10000750
$ws.Range("I13").Value = 10000750
$ws.Range("K13").Value = 10000750
$ws.Range("M13").Value = -10000606

$ws.Range("H61").Value = 3343.5715
$ws.Range("I61").Value = 3233.3333
$ws.Range("K61").Value = 3233.3333
$ws.Range("M61").Value = -3021.3333

$ws.Range("H63").Value = 1522.6818
$ws.Range("I63").Value = 1522.6818
$ws.Range("K63").Value = 1522.6818
$ws.Range("M63").Value = -836.6818000000001

$ws.Range("H66").Value = 1522.6818
$ws.Range("I66").Value = 1522.6818
$ws.Range("K66").Value = 7613.409000000001
$ws.Range("M66").Value = -4181.409000000001

$ws.Range("H132").Value = 3799.6667
$ws.Range("I132").Value = 1199.5
$ws.Range("K132").Value = 3598.5
$ws.Range("M132").Value = -1068.5

$ws.Range("H136").Value = 3343.5715
$ws.Range("I136").Value = 3233.3333
$ws.Range("K136").Value = 9699.999899999999
$ws.Range("M136").Value = -7149.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 540.375
$ws.Range("I16").Value = 408
$ws.Range("J16").Value = 937.5
$ws.Range("K16").Value = 408
$ws.Range("L16").Value = 937.5
$ws.Range("M16").Value = -238
$ws.Range("N16").Value = -1277.5

$ws.Range("H86").Value = 3551.2856
$ws.Range("I86").Value = 3134
$ws.Range("J86").Value = 3864.25
$ws.Range("K86").Value = 3134
$ws.Range("L86").Value = 3864.25
$ws.Range("M86").Value = -2011
$ws.Range("N86").Value = -6110.25

$ws.Range("H89").Value = 3551.2856
$ws.Range("I89").Value = 3134
$ws.Range("J89").Value = 3864.25
$ws.Range("K89").Value = 15670
$ws.Range("L89").Value = 19321.25
$ws.Range("M89").Value = -10054
$ws.Range("N89").Value = -30553.25

$ws.Range("H133").Value = 100780
$ws.Range("J133").Value = 100780
$ws.Range("L133").Value = 100780
$ws.Range("N133").Value = -110900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1448.8823
$ws.Range("I7").Value = 1340.5
$ws.Range("K7").Value = 1340.5
$ws.Range("M7").Value = -1227.5

$ws.Range("H16").Value = 1384.25
$ws.Range("I16").Value = 680
$ws.Range("K16").Value = 680
$ws.Range("M16").Value = -393

$ws.Range("H22").Value = 1096
$ws.Range("I22").Value = 1294.6666
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 1294.6666
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -944.6666
$ws.Range("N22").Value = -1200

$ws.Range("H58").Value = 999.5
$ws.Range("I58").Value = 999.5
$ws.Range("K58").Value = 999.5
$ws.Range("M58").Value = -796.5

$ws.Range("H97").Value = 16000
$ws.Range("J97").Value = 16000
$ws.Range("L97").Value = 16000
$ws.Range("N97").Value = -17982

$ws.Range("H99").Value = 6479.4375
$ws.Range("I99").Value = 6054.8335
$ws.Range("K99").Value = 6054.8335
$ws.Range("M99").Value = -4556.8335

$ws.Range("H113").Value = 1384.25
$ws.Range("I113").Value = 680
$ws.Range("K113").Value = 680
$ws.Range("M113").Value = 1490

$ws.Range("H126").Value = 6479.4375
$ws.Range("I126").Value = 6054.8335
$ws.Range("K126").Value = 18164.5005
$ws.Range("M126").Value = -15694.5005

$ws.Range("H132").Value = 6645.0527
$ws.Range("J132").Value = 13428.571
$ws.Range("L132").Value = 40285.713
$ws.Range("N132").Value = -45345.713

$ws.Range("H136").Value = 999.5
$ws.Range("I136").Value = 999.5
$ws.Range("K136").Value = 2998.5
$ws.Range("M136").Value = -448.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 717.4545000000001
$ws.Range("I92").Value = 498.66666
$ws.Range("J92").Value = 980
$ws.Range("K92").Value = 1495.99998
$ws.Range("L92").Value = 2940
$ws.Range("M92").Value = -247.9999800000001
$ws.Range("N92").Value = -5436

$ws.Range("H98").Value = 599.375
$ws.Range("I98").Value = 561.75
$ws.Range("J98").Value = 637
$ws.Range("K98").Value = 1685.25
$ws.Range("L98").Value = 1911
$ws.Range("M98").Value = -187.25
$ws.Range("N98").Value = -4907

$ws.Range("H129").Value = 3998.5715
$ws.Range("J129").Value = 4198
$ws.Range("L129").Value = 12594
$ws.Range("N129").Value = -22594

$ws.Range("H139").Value = 4558
$ws.Range("I139").Value = 2800
$ws.Range("K139").Value = 8400
$ws.Range("M139").Value = -3260

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 39435
$ws.Range("I20").Value = 4900
$ws.Range("J20").Value = 42888.5
$ws.Range("K20").Value = 4900
$ws.Range("L20").Value = 42888.5
$ws.Range("N20").Value = -43378.5
$ws.Range("M20").Value = -4655

$ws.Range("H113").Value = 2132.3333
$ws.Range("I113").Value = 1455.8572
$ws.Range("K113").Value = 1455.8572
$ws.Range("M113").Value = 714.1428000000001

$ws.Range("H122").Value = 127685.25
$ws.Range("I122").Value = 145356.28
$ws.Range("J122").Value = 3988
$ws.Range("K122").Value = 436068.84
$ws.Range("L122").Value = 11964
$ws.Range("M122").Value = -433618.84
$ws.Range("N122").Value = -16864

$ws.Range("H126").Value = 333335330
$ws.Range("I126").Value = 500002000
$ws.Range("K126").Value = 1500006000
$ws.Range("M126").Value = -1500003530

$ws.Range("H132").Value = 3221.2222
$ws.Range("I132").Value = 2498.6667
$ws.Range("K132").Value = 7496.000100000001
$ws.Range("M132").Value = -4966.000100000001

$ws.Range("H135").Value = 75780
$ws.Range("J135").Value = 75780
$ws.Range("L135").Value = 75780
$ws.Range("N135").Value = -85920

$ws.Range("H140").Value = 52500
$ws.Range("J140").Value = 52500
$ws.Range("L140").Value = 52500
$ws.Range("N140").Value = -62860

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 287166.84
$ws.Range("I10").Value = 2000000
$ws.Range("J10").Value = 1694.6666
$ws.Range("K10").Value = 2000000
$ws.Range("L10").Value = 1694.6666
$ws.Range("M10").Value = -1999860
$ws.Range("N10").Value = -1974.6666

$ws.Range("H19").Value = 4266.5
$ws.Range("I19").Value = 557.75
$ws.Range("J19").Value = 7975.25
$ws.Range("K19").Value = 557.75
$ws.Range("L19").Value = 7975.25
$ws.Range("M19").Value = -387.75
$ws.Range("N19").Value = -8315.25

$ws.Range("H22").Value = 4232.6665
$ws.Range("J22").Value = 4232.6665
$ws.Range("L22").Value = 4232.6665
$ws.Range("N22").Value = -4822.6665

$ws.Range("H27").Value = 4232.6665
$ws.Range("J27").Value = 4232.6665
$ws.Range("L27").Value = 4232.6665
$ws.Range("N27").Value = -4446.6665

$ws.Range("H61").Value = 2131.3572
$ws.Range("I61").Value = 1601.75
$ws.Range("K61").Value = 1601.75
$ws.Range("M61").Value = -1399.75

$ws.Range("H93").Value = 900.3333
$ws.Range("I93").Value = 900.4286
$ws.Range("J93").Value = 900
$ws.Range("K93").Value = 900.4286
$ws.Range("L93").Value = 900
$ws.Range("M93").Value = 347.5714
$ws.Range("N93").Value = -3396

$ws.Range("H113").Value = 2131.3572
$ws.Range("I113").Value = 1601.75
$ws.Range("K113").Value = 1601.75
$ws.Range("M113").Value = 568.25

$ws.Range("H132").Value = 7927.2856
$ws.Range("I132").Value = 7098.2
$ws.Range("K132").Value = 21294.6
$ws.Range("M132").Value = -18764.6

$ws.Range("H136").Value = 5400
$ws.Range("I136").Value = 5000
$ws.Range("K136").Value = 15000
$ws.Range("M136").Value = -12450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 31666.666
$ws.Range("J74").Value = 31666.666
$ws.Range("L74").Value = 31666.666
$ws.Range("N74").Value = -33538.666

$ws.Range("H77").Value = 31666.666
$ws.Range("J77").Value = 31666.666
$ws.Range("L77").Value = 94999.99800000001
$ws.Range("N77").Value = -104359.998

$ws.Range("H94").Value = 33330
$ws.Range("J94").Value = 33330
$ws.Range("L94").Value = 33330
$ws.Range("N94").Value = -35132

$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H136").Value = 3641.1
$ws.Range("J136").Value = 4665.6665
$ws.Range("L136").Value = 13996.9995
$ws.Range("N136").Value = -19096.9995

$ws.Range("H140").Value = 89429
$ws.Range("J140").Value = 89429
$ws.Range("L140").Value = 89429
$ws.Range("N140").Value = -99789

$ws.Range("H141").Value = 98000
$ws.Range("J141").Value = 98000
$ws.Range("L141").Value = 98000
$ws.Range("N141").Value = -108360
